# Update Name of Algo
# Applies corrected RandomForest imputation results to column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = -7.361199999999996
    6   = -7.7225
    7   = -7.774999999999999
    8   = -8.6751
    16  = -8.511700000000003
    20  = -8.015299999999991
    21  = -7.894099999999994
    28  = -8.066599999999999
    29  = -7.334100000000002
    30  = -7.547600000000002
    32  = -9.180099999999999
    40  = -7.951599999999996
    46  = -8.080199999999996
    51  = -7.674100000000001
    52  = -7.627999999999997
    57  = -8.252199999999995
    59  = -8.287299999999998
    62  = -9.086299999999989
    66  = -7.172800000000003
    73  = -7.456399999999993
    74  = -8.171900000000004
    77  = -5.860199999999997
    92  = -6.169400000000004
    100 = -8.302400000000004
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 4).Value = $updates[$row]
}
